$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")

# ALC row 40
$ws_ALC.Cells.Item(40, 8).Value = 3985.7144
$ws_ALC.Cells.Item(40, 10).Value = 4150
$ws_ALC.Cells.Item(40, 12).Value = 4150
$ws_ALC.Cells.Item(40, 14).Value = -4500

# ALC row 42
$ws_ALC.Cells.Item(42, 8).Value = 1279.6
$ws_ALC.Cells.Item(42, 10).Value = 1999.5
$ws_ALC.Cells.Item(42, 12).Value = 5998.5
$ws_ALC.Cells.Item(42, 14).Value = -6458.5

# ALC row 44
$ws_ALC.Cells.Item(44, 8).Value = 150000.5
$ws_ALC.Cells.Item(44, 10).Value = 150000.5
$ws_ALC.Cells.Item(44, 12).Value = 150000.5
$ws_ALC.Cells.Item(44, 14).Value = -150924.5

# ALC row 62
$ws_ALC.Cells.Item(62, 8).Value = 6108.4287
$ws_ALC.Cells.Item(62, 9).Value = 4938.25
$ws_ALC.Cells.Item(62, 11).Value = 4938.25
$ws_ALC.Cells.Item(62, 13).Value = -4314.25

# ALC row 65
$ws_ALC.Cells.Item(65, 8).Value = 6108.4287
$ws_ALC.Cells.Item(65, 9).Value = 4938.25
$ws_ALC.Cells.Item(65, 11).Value = 24691.25
$ws_ALC.Cells.Item(65, 13).Value = -21571.25

# ALC row 76
$ws_ALC.Cells.Item(76, 8).Value = 4996.625
$ws_ALC.Cells.Item(76, 9).Value = 4996.625
$ws_ALC.Cells.Item(76, 11).Value = 4996.625
$ws_ALC.Cells.Item(76, 13).Value = -4681.625

# ALC row 79
$ws_ALC.Cells.Item(79, 8).Value = 4996.625
$ws_ALC.Cells.Item(79, 9).Value = 4996.625
$ws_ALC.Cells.Item(79, 11).Value = 4996.625
$ws_ALC.Cells.Item(79, 13).Value = -3904.625

# ALC row 92
$ws_ALC.Cells.Item(92, 8).Value = 681.8
$ws_ALC.Cells.Item(92, 9).Value = 857.6087
$ws_ALC.Cells.Item(92, 10).Value = 344.83334
$ws_ALC.Cells.Item(92, 11).Value = 857.6087
$ws_ALC.Cells.Item(92, 12).Value = 344.83334
$ws_ALC.Cells.Item(92, 13).Value = 390.3913
$ws_ALC.Cells.Item(92, 14).Value = -2840.83334

# ALC row 112
$ws_ALC.Cells.Item(112, 8).Value = 779617.9
$ws_ALC.Cells.Item(112, 10).Value = 909187.9399999999
$ws_ALC.Cells.Item(112, 12).Value = 2727563.82
$ws_ALC.Cells.Item(112, 14).Value = -2729779.82

# ARM row 2
$ws_ARM.Cells.Item(2, 8).Value = 2415.3684
$ws_ARM.Cells.Item(2, 9).Value = 1847.6428
$ws_ARM.Cells.Item(2, 11).Value = 1847.6428
$ws_ARM.Cells.Item(2, 13).Value = -1734.6428

# ARM row 4
$ws_ARM.Cells.Item(4, 8).Value = 449.92307
$ws_ARM.Cells.Item(4, 9).Value = 325
$ws_ARM.Cells.Item(4, 10).Value = 866.3333
$ws_ARM.Cells.Item(4, 11).Value = 325
$ws_ARM.Cells.Item(4, 12).Value = 866.3333
$ws_ARM.Cells.Item(4, 13).Value = -209
$ws_ARM.Cells.Item(4, 14).Value = -1098.3333

# ARM row 5
$ws_ARM.Cells.Item(5, 8).Value = 144.5
$ws_ARM.Cells.Item(5, 9).Value = 144.5
$ws_ARM.Cells.Item(5, 11).Value = 144.5
$ws_ARM.Cells.Item(5, 13).Value = -32.5

# ARM row 41
$ws_ARM.Cells.Item(41, 8).Value = 18145.75
$ws_ARM.Cells.Item(41, 9).Value = 5700
$ws_ARM.Cells.Item(41, 10).Value = 38888.668
$ws_ARM.Cells.Item(41, 11).Value = 5700
$ws_ARM.Cells.Item(41, 12).Value = 38888.668
$ws_ARM.Cells.Item(41, 13).Value = -5286
$ws_ARM.Cells.Item(41, 14).Value = -39716.668

# ARM row 74
$ws_ARM.Cells.Item(74, 8).Value = 1960.7084
$ws_ARM.Cells.Item(74, 9).Value = 1602.1818
$ws_ARM.Cells.Item(74, 11).Value = 1602.1818
$ws_ARM.Cells.Item(74, 13).Value = -728.1818000000001

# ARM row 77
$ws_ARM.Cells.Item(77, 8).Value = 1960.7084
$ws_ARM.Cells.Item(77, 9).Value = 1602.1818
$ws_ARM.Cells.Item(77, 11).Value = 8010.909000000001
$ws_ARM.Cells.Item(77, 13).Value = -3642.909000000001

# ARM row 116
$ws_ARM.Cells.Item(116, 8).Value = 2415.3684
$ws_ARM.Cells.Item(116, 9).Value = 1847.6428
$ws_ARM.Cells.Item(116, 11).Value = 1847.6428
$ws_ARM.Cells.Item(116, 13).Value = 446.3571999999999

# BSM row 3
$ws_BSM.Cells.Item(3, 8).Value = 2415.3684
$ws_BSM.Cells.Item(3, 9).Value = 1847.6428
$ws_BSM.Cells.Item(3, 11).Value = 1847.6428
$ws_BSM.Cells.Item(3, 13).Value = -1733.6428

# BSM row 4
$ws_BSM.Cells.Item(4, 8).Value = 144.5
$ws_BSM.Cells.Item(4, 9).Value = 144.5
$ws_BSM.Cells.Item(4, 11).Value = 144.5
$ws_BSM.Cells.Item(4, 13).Value = -29.5

# BSM row 69
$ws_BSM.Cells.Item(69, 8).Value = 60000
$ws_BSM.Cells.Item(69, 10).Value = 60000
$ws_BSM.Cells.Item(69, 12).Value = 60000
$ws_BSM.Cells.Item(69, 14).Value = -61622

# BSM row 72
$ws_BSM.Cells.Item(72, 8).Value = 60000
$ws_BSM.Cells.Item(72, 10).Value = 60000
$ws_BSM.Cells.Item(72, 12).Value = 180000
$ws_BSM.Cells.Item(72, 14).Value = -188112

# CRP row 2
$ws_CRP.Cells.Item(2, 8).Value = 1924.2222
$ws_CRP.Cells.Item(2, 9).Value = 2463
$ws_CRP.Cells.Item(2, 10).Value = 846.6667
$ws_CRP.Cells.Item(2, 11).Value = 2463
$ws_CRP.Cells.Item(2, 12).Value = 846.6667
$ws_CRP.Cells.Item(2, 13).Value = -2350
$ws_CRP.Cells.Item(2, 14).Value = -1072.6667

# CRP row 3
$ws_CRP.Cells.Item(3, 8).Value = 640.3333
$ws_CRP.Cells.Item(3, 9).Value = 582.875
$ws_CRP.Cells.Item(3, 10).Value = 1100
$ws_CRP.Cells.Item(3, 11).Value = 582.875
$ws_CRP.Cells.Item(3, 12).Value = 1100
$ws_CRP.Cells.Item(3, 13).Value = -469.875
$ws_CRP.Cells.Item(3, 14).Value = -1326

# CRP row 7
$ws_CRP.Cells.Item(7, 8).Value = 1097.6875
$ws_CRP.Cells.Item(7, 9).Value = 254.64285
$ws_CRP.Cells.Item(7, 11).Value = 254.64285
$ws_CRP.Cells.Item(7, 13).Value = -141.64285

# CRP row 86
$ws_CRP.Cells.Item(86, 8).Value = 3678.4285
$ws_CRP.Cells.Item(86, 9).Value = 3416.6667
$ws_CRP.Cells.Item(86, 11).Value = 3416.6667
$ws_CRP.Cells.Item(86, 13).Value = -2293.6667

# CRP row 89
$ws_CRP.Cells.Item(89, 8).Value = 3678.4285
$ws_CRP.Cells.Item(89, 9).Value = 3416.6667
$ws_CRP.Cells.Item(89, 11).Value = 17083.3335
$ws_CRP.Cells.Item(89, 13).Value = -11467.3335

# CUL row 2
$ws_CUL.Cells.Item(2, 8).Value = 85.333336
$ws_CUL.Cells.Item(2, 9).Value = 59.5
$ws_CUL.Cells.Item(2, 10).Value = 98.25
$ws_CUL.Cells.Item(2, 11).Value = 357
$ws_CUL.Cells.Item(2, 12).Value = 589.5
$ws_CUL.Cells.Item(2, 13).Value = -244
$ws_CUL.Cells.Item(2, 14).Value = -815.5

# CUL row 92
$ws_CUL.Cells.Item(92, 8).Value = 687.1429000000001
$ws_CUL.Cells.Item(92, 9).Value = 606.75
$ws_CUL.Cells.Item(92, 10).Value = 794.3333
$ws_CUL.Cells.Item(92, 11).Value = 1820.25
$ws_CUL.Cells.Item(92, 12).Value = 2382.9999
$ws_CUL.Cells.Item(92, 13).Value = -572.25
$ws_CUL.Cells.Item(92, 14).Value = -4878.9999

# CUL row 93
$ws_CUL.Cells.Item(93, 8).Value = 14577.4
$ws_CUL.Cells.Item(93, 9).Value = 6999.5
$ws_CUL.Cells.Item(93, 10).Value = 19629.334
$ws_CUL.Cells.Item(93, 11).Value = 20998.5
$ws_CUL.Cells.Item(93, 12).Value = 58888.00199999999
$ws_CUL.Cells.Item(93, 13).Value = -19126.5
$ws_CUL.Cells.Item(93, 14).Value = -62632.00199999999

# CUL row 131
$ws_CUL.Cells.Item(131, 8).Value = 800841
$ws_CUL.Cells.Item(131, 9).Value = 126060.75
$ws_CUL.Cells.Item(131, 10).Value = 1138231.1
$ws_CUL.Cells.Item(131, 11).Value = 378182.25
$ws_CUL.Cells.Item(131, 12).Value = 3414693.3
$ws_CUL.Cells.Item(131, 13).Value = -373142.25
$ws_CUL.Cells.Item(131, 14).Value = -3424773.3

# GSM row 113
$ws_GSM.Cells.Item(113, 8).Value = 3880.1
$ws_GSM.Cells.Item(113, 9).Value = 3721.875
$ws_GSM.Cells.Item(113, 10).Value = 4513
$ws_GSM.Cells.Item(113, 11).Value = 3721.875
$ws_GSM.Cells.Item(113, 12).Value = 4513
$ws_GSM.Cells.Item(113, 13).Value = -1551.875
$ws_GSM.Cells.Item(113, 14).Value = -8853

# GSM row 126
$ws_GSM.Cells.Item(126, 8).Value = 8151.2104
$ws_GSM.Cells.Item(126, 9).Value = 6135.1113
$ws_GSM.Cells.Item(126, 10).Value = 9965.700000000001
$ws_GSM.Cells.Item(126, 11).Value = 18405.3339
$ws_GSM.Cells.Item(126, 12).Value = 29897.1
$ws_GSM.Cells.Item(126, 13).Value = -15935.3339
$ws_GSM.Cells.Item(126, 14).Value = -34837.10000000001

# LTW row 94
$ws_LTW.Cells.Item(94, 8).Value = 0
$ws_LTW.Cells.Item(94, 10).Value = 0
$ws_LTW.Cells.Item(94, 12).Value = 0
$ws_LTW.Cells.Item(94, 14).ClearContents()

# LTW row 100
$ws_LTW.Cells.Item(100, 8).Value = 335091.44
$ws_LTW.Cells.Item(100, 9).Value = 376746.62
$ws_LTW.Cells.Item(100, 10).Value = 1850
$ws_LTW.Cells.Item(100, 11).Value = 376746.62
$ws_LTW.Cells.Item(100, 12).Value = 1850
$ws_LTW.Cells.Item(100, 13).Value = -376205.62
$ws_LTW.Cells.Item(100, 14).Value = -2932
